$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update C9: 0 -> 1 (this drives the recalculated averages in C15/I15)
$ws.Range("C9").Value = 1

# Clear J9, which held the note "Need shortBy method"
$ws.Range("J9").ClearContents()

# Update the active selection to match the author's final cursor position
$ws.Range("E9").Select()
